$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Insert 4 new rows (10-13) into the "Issues" sheet, pushing existing rows down ---
$ws1.Rows.Item(10).Resize(4).Insert()

$ws1.Range("A10").Value = "Model arbitarily started in 1974 (50 historical years)"
$ws1.Range("B10").Value = "What is a suitable time to start fishing?"
$ws1.Rows.Item(10).RowHeight = 29.15

$ws1.Range("A11").Value = "Equilibrium catches assumed to be negligible before 1974 (again arbitrarily) of those observed"
$ws1.Range("B11").Value = "C_eq = 0.00001 * mean historical catch but can bring forward the initial model year and specify a differing C_eq"
$ws1.Rows.Item(11).RowHeight = 43.75

$ws1.Range("A12").Value = "OM fleet and survey structure: "
$ws1.Range("B12").Value = "2 fleets (rec / commercial), 4 surveys (rec CPUE, commercial CPUE, historial length/age comp, rec survey length/age comp)"
$ws1.Rows.Item(12).RowHeight = 29.15

$ws1.Range("A13").Value = "Fleet seleectivities are not informed by length / age data"
$ws1.Range("B13").Value = "Currently specified - but can we assume that the rec survey reflects the recreational fishery? If so then we can move around the data to make the rec selectivity / CPUE informated by the same length/age observations"
$ws1.Rows.Item(13).RowHeight = 43.75

# --- Rename the "Commercial fishery exploitation is unknown" issue (now row 24 after the insert) ---
$ws1.Range("A24").Value = "Commercial fishery selectivity / retention is unknown"
$ws1.Range("B24").Value = "For now I'm assuming it follows the survey comp data"

# --- Update the selection / view on the Issues sheet ---
$ws1.Range("O12").Select()

# --- Add the new "Todo" sheet after "Issues" ---
$todo = $wb.Worksheets.Add($null, $ws1)
$todo.Name = "Todo"

$todo.Columns.Item(1).ColumnWidth = 22.53515625
$todo.Columns.Item(2).ColumnWidth = 27.4609375

$todo.Range("A1").Value = "High"
$todo.Range("B1").Value = "Meet to discuss straw-dog fits"

$todo.Range("A2").Value = "Medium"
$todo.Range("B2").Value = "Data weighting profiling"

$todo.Range("A3").Value = "Medium"
$todo.Range("B3").Value = "Parameter profiling"

$todo.Range("B22").Select()
